# Auto update price data 2025-11-28
#
# The daily price-tracking sheet gets a new row inserted right under the
# header for the latest trading day. All previously recorded rows shift
# down by one (their dates/values are untouched, they just move from row N
# to row N+1), and the newest date/price observation is written into the
# freshly inserted row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the first data row (row 2); rows 2.. shift to 3..
$ws.Rows.Item(2).Insert()

# The inserted row picks up formatting copied from its neighbour, so reset
# it back to the plain "Normal" style used by the rest of the data rows.
$ws.Range("A2:D2").Style = "Normal"

# Force column A to remain plain text (matching the other date cells)
# instead of being auto-recognized as a date value.
$ws.Cells.Item(2, 1).NumberFormat = "@"

# Latest day's data
$ws.Cells.Item(2, 1).Value = "2025-11-28"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
